$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 285, shifting existing rows 285.. down by one.
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the new record's data.
$ws.Cells.Item(285, 1).Value = 5
$ws.Cells.Item(285, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(285, 3).Value = "Maule"
$ws.Cells.Item(285, 4).Value = 44694
$ws.Cells.Item(285, 5).Value = 7
$ws.Cells.Item(285, 6).Value = 100114013
$ws.Cells.Item(285, 7).Value = "Zanahoria"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 500
$ws.Cells.Item(285, 11).Value = 6000
$ws.Cells.Item(285, 12).Value = 6000
$ws.Cells.Item(285, 13).Value = 6000
$ws.Cells.Item(285, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(285, 15).Value = "Región de Ñuble"
$ws.Cells.Item(285, 16).Value = 300
$ws.Cells.Item(285, 17).Value = 20
$ws.Cells.Item(285, 18).Value = "Hortaliza"
